$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 2815608.16
$ws.Range("C9").Value = 436776.99
$ws.Range("D9").Value = 3252385.15
$ws.Range("E9").Value = 13.42943624004678
$ws.Range("F9").Value = 86.57056375995322
$ws.Range("G9").Value = -57.78768178697559
$ws.Range("H9").Value = -49.15400847553775
$ws.Range("I9").Value = -50.51326976586423
$ws.Range("J9").Value = 27908
$ws.Range("K9").Value = 1186
$ws.Range("L9").Value = 29094
